# Auto-generated Excel COM-interop script to update cryptos.xlsx data
# Applies value changes for rows 2-51 (B=Coin, C=Link, D=Price, E=Volume(1h))
# Matches the "Updated cryptos list ... with GitHub Actions" commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "69.216.95"
$ws.Range("E2").Value = "  +1.32%  "
# Row 3
$ws.Range("D3").Value = "3.775.80"
$ws.Range("E3").Value = "  -0.52%  "
# Row 4
$ws.Range("E4").Value = "  +0.30%  "
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "630.12"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.69%  "
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "166.02"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.67%  "
# Row 7
$ws.Range("D7").Value = "3.774.53"
$ws.Range("E7").Value = "  -0.50%  "
# Row 8
$ws.Range("E8").Value = "  -0.05%  "
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.521"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.79%  "
# Row 10
$ws.Range("E10").Value = "  -0.30%  "
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.460"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.31%  "
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.82"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.44%  "
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000243"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.67%  "
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.81"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.63%  "
# Row 15
$ws.Range("D15").Value = "4.409.43"
$ws.Range("E15").Value = "  -0.53%  "
# Row 16
$ws.Range("D16").Value = "3.790.43"
$ws.Range("E16").Value = "  +1.07%  "
# Row 17
$ws.Range("D17").Value = "69.219.14"
$ws.Range("E17").Value = "  +1.37%  "
# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "17.62"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.42%  "
# Row 19
$ws.Range("E19").Value = "  +0.45%  "
# Row 20
$ws.Range("E20").Value = "  -0.52%  "
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "464.05"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.51%  "
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.54"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.34%  "
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.707"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.21%  "
# Row 24
$ws.Range("B24").Value = "PEPE"
$ws.Range("C24").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0000145"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.38%  "
# Row 25
$ws.Range("B25").Value = "Litecoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "82.90"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.70%  "
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.04"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.52%  "
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.14"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.85%  "
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.14"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.67%  "
# Row 29
$ws.Range("E29").Value = "  +0.02%  "
# Row 30
$ws.Range("D30").Value = "3.925.88"
$ws.Range("E30").Value = "  -0.55%  "
# Row 31
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.68"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.11%  "
# Row 32
$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.28"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.01%  "
# Row 33
$ws.Range("E33").Value = "  -1.85%  "
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "28.55"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.58%  "
# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.170"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +14.31%  "
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.999"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.02%  "
# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "8.99"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.29%  "
# Row 38
$ws.Range("D38").Value = "3.728.99"
$ws.Range("E38").Value = "  -0.38%  "
# Row 39
$ws.Range("E39").Value = "  +0.20%  "
# Row 40
$ws.Range("E40").Value = "  +4.21%  "
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.80"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.24%  "
# Row 42
$ws.Range("B42").Value = "FirstDigitalUSD"
$ws.Range("C42").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.00"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.10%  "
# Row 43
$ws.Range("B43").Value = "Mantle"
$ws.Range("C43").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.960"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.61%  "
# Row 44
$ws.Range("E44").Value = "  -0.01%  "
# Row 45
$ws.Range("E45").Value = "  +3.26%  "
# Row 46
$ws.Range("B46").Value = "ONDO"
$ws.Range("C46").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.43"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.84%  "
# Row 47
$ws.Range("B47").Value = "Stacks"
$ws.Range("C47").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.95"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.97%  "
# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "43.01"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.29%  "
# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.295"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.09%  "
# Row 50
$ws.Range("E50").Value = "  -0.72%  "
# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.37"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.10%  "

Write-Host "Applied 104 cell updates"
